$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 675 (the post about "「あす日が昇るだろう。私は朝が大好きだ」"),
# shifting all subsequent rows up by one.
$ws.Rows.Item(675).Delete()
